$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted before the existing row 997, pushing the
# former rows 997-1071 down to 998-1072 (dimension grows from R1071 to R1072).
$ws.Rows.Item(997).Insert()

$ws.Range("A997").Value() = 6
$ws.Range("B997").Value() = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C997").Value() = 'Metropolitana'
$ws.Range("D997").Value() = 45013
$ws.Range("E997").Value() = 13
$ws.Range("F997").Value() = 100112003
$ws.Range("G997").Value() = 'Ajo'
$ws.Range("H997").Value() = 'Chino'
$ws.Range("I997").Value() = 'Primera'
$ws.Range("J997").Value() = 1400
$ws.Range("K997").Value() = 13500
$ws.Range("L997").Value() = 14000
$ws.Range("M997").Value() = 13679
$ws.Range("N997").Value() = '$/caja 10 kilos'
$ws.Range("O997").Value() = 'China'
$ws.Range("P997").Value() = 1368
$ws.Range("Q997").Value() = 10
$ws.Range("R997").Value() = 'Hortaliza'
